# Update "Intervention content and delivery.xlsx"
# Commit: Updates LSR numbers in spreadsheet
#
# The "LSR no." column (P) used to hold strings like "LSR 3" or
# "LSR 1; LSR 2" - these are rewritten as plain numbers ("3") or
# semi-colon separated lists without the "LSR " prefix ("1; 2").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column P = 16 ("LSR no.")
# Single LSR values become plain numbers.
$ws.Cells.Item(2,16).Value  = 3
$ws.Cells.Item(3,16).Value  = 3
$ws.Cells.Item(4,16).Value  = 3
$ws.Cells.Item(5,16).Value  = 3

# Row 7 before row 6 so that the new shared-string combo values are
# appended to the shared string table in the same order as the source.
$ws.Cells.Item(7,16).Value  = "2; 3"
$ws.Cells.Item(6,16).Value  = "1; 3"

$ws.Cells.Item(8,16).Value  = 3
$ws.Cells.Item(9,16).Value  = 3
$ws.Cells.Item(10,16).Value = 3
$ws.Cells.Item(11,16).Value = 3
$ws.Cells.Item(12,16).Value = 1

$ws.Cells.Item(13,16).Value = "1; 2; 3"
$ws.Cells.Item(14,16).Value = "1; 2; 3"

$ws.Cells.Item(15,16).Value = 3
$ws.Cells.Item(16,16).Value = 3
$ws.Cells.Item(17,16).Value = 3
$ws.Cells.Item(18,16).Value = 3
$ws.Cells.Item(19,16).Value = 3
$ws.Cells.Item(20,16).Value = 3
$ws.Cells.Item(21,16).Value = 3

$ws.Cells.Item(22,16).Value = "1; 2; 3"

$ws.Cells.Item(23,16).Value = 3
$ws.Cells.Item(24,16).Value = 3
$ws.Cells.Item(25,16).Value = 3
$ws.Cells.Item(26,16).Value = 3
$ws.Cells.Item(27,16).Value = 1
$ws.Cells.Item(28,16).Value = 1

$ws.Cells.Item(29,16).Value = "2; 3"
$ws.Cells.Item(30,16).Value = "1; 2"
$ws.Cells.Item(31,16).Value = "1; 2"

$ws.Cells.Item(32,16).Value = 1

$ws.Cells.Item(33,16).Value = 2
$ws.Cells.Item(34,16).Value = 2
$ws.Cells.Item(35,16).Value = 2
$ws.Cells.Item(36,16).Value = 2
$ws.Cells.Item(37,16).Value = 2
$ws.Cells.Item(38,16).Value = 2
$ws.Cells.Item(39,16).Value = 2

# Rows 34-37 (columns A:D) also lost their (invisible, fillId=0) fill
# formatting flag in this save - clear the interior fill so the cell
# styles collapse back to the plain bordered/top-aligned styles used
# elsewhere in the sheet.
$ws.Range("A34:D37").Interior.Pattern = -4142  # xlNone

# A36/A37 were empty placeholder cells that only existed to carry that
# now-removed formatting; drop them entirely.
$ws.Range("A36:A37").Clear()

# Leave the selection where the author's last edit was.
$ws.Range("T14").Select()
